$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($count)
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$new.Name = "magapoke_2026-02-25"

$new.Cells.Item(1,1).Value = "rank"
$new.Cells.Item(1,2).Value = "title"

$titles = @(
  "ブルーロック",
  "東京卍リベンジャーズ",
  "ギルティサークル",
  "ベイビーステップ",
  "島耕作",
  "君が僕らを悪魔と呼んだ頃",
  "イレギュラーズ",
  "黄昏町プリズナーズ",
  "信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！",
  "十字架のろくにん",
  "愛妻の裏アカ",
  "魔女と傭兵",
  "転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～",
  "黒猫と魔女の教室",
  "ハードワーカー中田",
  "ガチアクタ",
  "転生したら第七王子だったので、気ままに魔術を極めます",
  "となりの黒川さん",
  "南海トラフ巨大地震",
  "魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～",
  "デッドアカウント",
  "ドラハチ",
  "WIND BREAKER",
  "【爆アド】生まれた直後から最強悪霊と脳内バトルしてたら魔力量が測定可能域を超えてました〜悪憑の子の謙虚な覇道〜",
  "異世界ウォーキング",
  "K-9~警視庁公安部公安第9課異能対策係~",
  "ひゃくえむ。",
  "限界集落を脱村した錬金術士、都会で`"最強`"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～",
  "蒼く染めろ",
  "味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す",
  "さわらないで小手指くん",
  "幼馴染とはラブコメにならない",
  "せいぶつ部の田辺くん",
  "グラぱらっ！",
  "追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜",
  "辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜",
  "ともだちづくり",
  "ハンドレッドノート－アグリーダック－",
  "ハナバス　苔石花江のバスケ論",
  "屋根の下のアルテミス",
  "アルキメデスの大戦",
  "お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！",
  "Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。",
  "FAIRY TAIL 100 YEARS QUEST",
  "食糧人類-Starving Anonymous-",
  "普通の本はありません！",
  "降り積もれ孤独な死よ",
  "おやすみ ふみさん",
  "異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～",
  "不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～",
  "皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～",
  "ペンの夢に紅をさす",
  "念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！",
  "ジュミドロ",
  "なれの果ての僕ら",
  "いじめるヤバイ奴",
  "Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜",
  "可愛いだけじゃない式守さん",
  "アオバノバスケ",
  "時々ボソッとロシア語でデレる隣のアーリャさん",
  "君が監督！",
  "東京卍リベンジャーズ～場地圭介からの手紙～",
  "恋ニ非ズ",
  "MYS",
  "不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～",
  "田んぼで拾った女騎士、田舎で俺の嫁だと思われている",
  "剣帝学院の魔眼賢者",
  "白鳥運子は31画",
  "おくることば",
  "ストーカー行為がバレて人生終了男",
  "追放されなかった男　～二度目の人生は土下座から始まりました～",
  "死ぬほど君の処女が欲しい",
  "インフェクション",
  "シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～",
  "春くらり",
  "最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～",
  "デスティニーラバーズ",
  "GALAXIAS",
  "私をセンターにすると誓いますか？",
  "我間乱 ―修羅―",
  "はっちぽっちぱんち",
  "この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～",
  "鳴るさんだぁ",
  "劣等人の魔剣使い　スキルボードを駆使して最強に至る",
  "ヒロインは絶望しました。",
  "阿武ノーマル",
  "ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～",
  "This Man その顔を見た者には死を",
  "復讐の教科書",
  "イジらないで、長瀞さん",
  "ほねぬきごはん　～ウブで奥手な半キュバスにハートをください～",
  "ぼくたちのリメイク",
  "冰剣の魔術師が世界を統べる〜世界最強の魔術師である少年は、魔術学院に入学する〜",
  "ハンドレッドノート－高校生探偵 天命大地－",
  "「無能はいらない」と言われたから絶縁してやった　～最強の四天王に育てられた俺は、冒険者となり無双する～",
  "ザ・ファブル",
  "五輪の女神さま 〜なでしこ寮のメダルごはん〜",
  "魁の花巫女",
  "彼女、お借りします",
  "メダリスト"
)

for ($i = 0; $i -lt $titles.Length; $i++) {
  $row = $i + 2
  $new.Cells.Item($row, 1).Value = $i + 1
  $new.Cells.Item($row, 2).Value = $titles[$i]
}

Write-Output "Added sheet with $($titles.Length) rows"